$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row containing "atmospheric corrosion resistant steel" / "thép chống gỉ"
# (row 10). Remaining rows shift up, and the former last row (18) disappears.
$ws.Rows.Item(10).Delete()
